$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy donor_id (column A) into source (column C) for the newly-added
# non-directed donor (NDD) rows 246-263, matching the "source" to the
# donor's own id.
for ($r = 246; $r -le 263; $r++) {
    $donorId = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 3).Value = $donorId
}

# Reflect the user's on-screen selection after the edit (the last-touched
# range of newly-filled "source" cells).
$ws.Range("C246:C263").Select()
